# Samples for onRowSorted and onColumnSorted (Sort, Events)
# Adds 3 rows to the "Snippets" table describing the new sort-event snippets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by three rows so the table/autoFilter ref (and the
# sheet dimension) pick up rows 208:210 automatically.
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Duplicate the formatting of the last existing data row (206/207 use the
# highlighted banding style) down into the three new rows by copying it and
# inserting the copy - this carries the cell style along, unlike a plain
# value assignment.
$lastRow = $ws.Range("A207:D207")

$lastRow.Copy()
$ws.Range("A208:D208").Insert(-4121, 0)

$lastRow.Copy()
$ws.Range("A209:D209").Insert(-4121, 0)

$lastRow.Copy()
$ws.Range("A210:D210").Insert(-4121, 0)

# Row 208: Range.sort
$ws.Range("A208").Value = "Range"
$ws.Range("B208").Value = "sort"
$ws.Range("C208").Value = "excel-event-column-and-row-sort"
$ws.Range("D208").Value = "sortTopToBottom"

# Row 209: Worksheet.onRowSorted (write D before B so the shared-string
# table allocates "registerRowSortHandler" ahead of "onRowSorted").
$ws.Range("A209").Value = "Worksheet"
$ws.Range("D209").Value = "registerRowSortHandler"
$ws.Range("B209").Value = "onRowSorted"
$ws.Range("C209").Value = "excel-event-column-and-row-sort"

# Row 210: Worksheet.onColumnSorted
$ws.Range("A210").Value = "Worksheet"
$ws.Range("B210").Value = "onColumnSorted"
$ws.Range("C210").Value = "excel-event-column-and-row-sort"
$ws.Range("D210").Value = "registerColumnSortHandler"

# Scroll the frozen view down and land the selection on the last new cell,
# matching the author's final on-screen state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 183
$ws.Range("D210").Select()
